$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15-37 down to 16-38.
# Excel copies the formatting (incl. the date number format on column D)
# from the row above on insert, matching the target workbook's cell styles.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new Cilantro price record.
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = 44771
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = 100112040
$ws.Cells.Item(15, 7).Value = "Cilantro"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 15000
$ws.Cells.Item(15, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 417
$ws.Cells.Item(15, 17).Value = 36
$ws.Cells.Item(15, 18).Value = "Hortaliza"
